$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations": refresh the BRVM ranking table (rows 2-47) ---
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 2: NEI-CEDA CI
$ws1.Cells.Item(2,1).Value = "NEI-CEDA CI"
$ws1.Cells.Item(2,2).Value = 0
$ws1.Cells.Item(2,3).Value = 4
$ws1.Cells.Item(2,4).Value = 3680
$ws1.Cells.Item(2,5).Value = 945
$ws1.Cells.Item(2,6).Value = "🟡 Observer"
$ws1.Cells.Item(2,7).Value = "➖ Neutre"

# Row 3: BRVM - SERVICES PUBLICS
$ws1.Cells.Item(3,1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(3,2).Value = 0
$ws1.Cells.Item(3,3).Value = 8
$ws1.Cells.Item(3,4).Value = 3426.66
$ws1.Cells.Item(3,5).Value = 111.96
$ws1.Cells.Item(3,6).Value = "🟡 Observer"
$ws1.Cells.Item(3,7).Value = "➖ Neutre"

# Row 4: AIR LIQUIDE CI
$ws1.Cells.Item(4,1).Value = "AIR LIQUIDE CI"
$ws1.Cells.Item(4,2).Value = 0
$ws1.Cells.Item(4,3).Value = 4
$ws1.Cells.Item(4,4).Value = 2750
$ws1.Cells.Item(4,5).Value = 700
$ws1.Cells.Item(4,6).Value = "🟡 Observer"
$ws1.Cells.Item(4,7).Value = "➖ Neutre"

# Row 5: BRVM - AUTRES SECTEURS
$ws1.Cells.Item(5,1).Value = "BRVM - AUTRES SECTEURS"
$ws1.Cells.Item(5,2).Value = 0
$ws1.Cells.Item(5,3).Value = 4
$ws1.Cells.Item(5,4).Value = 2363.56
$ws1.Cells.Item(5,5).Value = 599.2
$ws1.Cells.Item(5,6).Value = "🟡 Observer"
$ws1.Cells.Item(5,7).Value = "➖ Neutre"

# Row 6: BRVM - DISTRIBUTION
$ws1.Cells.Item(6,1).Value = "BRVM - DISTRIBUTION"
$ws1.Cells.Item(6,2).Value = 0
$ws1.Cells.Item(6,3).Value = 4
$ws1.Cells.Item(6,4).Value = 2158.19
$ws1.Cells.Item(6,5).Value = 541.92
$ws1.Cells.Item(6,6).Value = "🟡 Observer"
$ws1.Cells.Item(6,7).Value = "➖ Neutre"

# Row 7: BRVM - TRANSPORT
$ws1.Cells.Item(7,1).Value = "BRVM - TRANSPORT"
$ws1.Cells.Item(7,2).Value = 0
$ws1.Cells.Item(7,3).Value = 4
$ws1.Cells.Item(7,4).Value = 1425.8
$ws1.Cells.Item(7,5).Value = 354.92
$ws1.Cells.Item(7,6).Value = "🟡 Observer"
$ws1.Cells.Item(7,7).Value = "➖ Neutre"

# Row 8: BRVM - AGRICULTURE
$ws1.Cells.Item(8,1).Value = "BRVM - AGRICULTURE"
$ws1.Cells.Item(8,2).Value = 0
$ws1.Cells.Item(8,3).Value = 4
$ws1.Cells.Item(8,4).Value = 1371.37
$ws1.Cells.Item(8,5).Value = 342.3
$ws1.Cells.Item(8,6).Value = "🟡 Observer"
$ws1.Cells.Item(8,7).Value = "➖ Neutre"

# Row 9: SUCRIVOIRE
$ws1.Cells.Item(9,1).Value = "SUCRIVOIRE"
$ws1.Cells.Item(9,2).Value = 0
$ws1.Cells.Item(9,3).Value = 1
$ws1.Cells.Item(9,4).Value = 990
$ws1.Cells.Item(9,5).Value = 990
$ws1.Cells.Item(9,6).Value = "🟡 Observer"
$ws1.Cells.Item(9,7).Value = "➖ Neutre"

# Row 10: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Cells.Item(10,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(10,2).Value = 0
$ws1.Cells.Item(10,3).Value = 4
$ws1.Cells.Item(10,4).Value = 748.48
$ws1.Cells.Item(10,5).Value = 189.13
$ws1.Cells.Item(10,6).Value = "🟡 Observer"
$ws1.Cells.Item(10,7).Value = "➖ Neutre"

# Row 11: BRVM - CONSOMMATION DE BASE         (**)
$ws1.Cells.Item(11,1).Value = "BRVM - CONSOMMATION DE BASE         (**)"
$ws1.Cells.Item(11,2).Value = 0
$ws1.Cells.Item(11,3).Value = 3
$ws1.Cells.Item(11,4).Value = 675.6
$ws1.Cells.Item(11,5).Value = 225.67
$ws1.Cells.Item(11,6).Value = "🟡 Observer"
$ws1.Cells.Item(11,7).Value = "➖ Neutre"

# Row 12: BRVM - FINANCES
$ws1.Cells.Item(12,1).Value = "BRVM - FINANCES"
$ws1.Cells.Item(12,2).Value = 0
$ws1.Cells.Item(12,3).Value = 4
$ws1.Cells.Item(12,4).Value = 588.66
$ws1.Cells.Item(12,5).Value = 148.41
$ws1.Cells.Item(12,6).Value = "🟡 Observer"
$ws1.Cells.Item(12,7).Value = "➖ Neutre"

# Row 13: BRVM - SERVICES FINANCIERS
$ws1.Cells.Item(13,1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(13,2).Value = 0
$ws1.Cells.Item(13,3).Value = 4
$ws1.Cells.Item(13,4).Value = 578.52
$ws1.Cells.Item(13,5).Value = 145.85
$ws1.Cells.Item(13,6).Value = "🟡 Observer"
$ws1.Cells.Item(13,7).Value = "➖ Neutre"

# Row 14: BRVM-PRESTIGE
$ws1.Cells.Item(14,1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(14,2).Value = 0
$ws1.Cells.Item(14,3).Value = 4
$ws1.Cells.Item(14,4).Value = 577.1
$ws1.Cells.Item(14,5).Value = 146
$ws1.Cells.Item(14,6).Value = "🟡 Observer"
$ws1.Cells.Item(14,7).Value = "➖ Neutre"

# Row 15: BRVM - INDUSTRIELS
$ws1.Cells.Item(15,1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(15,2).Value = 0
$ws1.Cells.Item(15,3).Value = 4
$ws1.Cells.Item(15,4).Value = 489.97
$ws1.Cells.Item(15,5).Value = 122.79
$ws1.Cells.Item(15,6).Value = "🟡 Observer"
$ws1.Cells.Item(15,7).Value = "➖ Neutre"

# Row 16: BRVM - ENERGIE
$ws1.Cells.Item(16,1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(16,2).Value = 0
$ws1.Cells.Item(16,3).Value = 4
$ws1.Cells.Item(16,4).Value = 454.56
$ws1.Cells.Item(16,5).Value = 113.76
$ws1.Cells.Item(16,6).Value = "🟡 Observer"
$ws1.Cells.Item(16,7).Value = "➖ Neutre"

# Row 17: BRVM-PRINCIPAL                    (**)
$ws1.Cells.Item(17,1).Value = "BRVM-PRINCIPAL                    (**)"
$ws1.Cells.Item(17,2).Value = 0
$ws1.Cells.Item(17,3).Value = 2
$ws1.Cells.Item(17,4).Value = 443.03
$ws1.Cells.Item(17,5).Value = 221.95
$ws1.Cells.Item(17,6).Value = "🟡 Observer"
$ws1.Cells.Item(17,7).Value = "➖ Neutre"

# Row 18: BRVM - TELECOMMUNICATIONS
$ws1.Cells.Item(18,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(18,2).Value = 0
$ws1.Cells.Item(18,3).Value = 4
$ws1.Cells.Item(18,4).Value = 387.98
$ws1.Cells.Item(18,5).Value = 97.25
$ws1.Cells.Item(18,6).Value = "🟡 Observer"
$ws1.Cells.Item(18,7).Value = "➖ Neutre"

# Row 19: BRVM - INDUSTRIE                 (**)
$ws1.Cells.Item(19,1).Value = "BRVM - INDUSTRIE                 (**)"
$ws1.Cells.Item(19,2).Value = 0
$ws1.Cells.Item(19,3).Value = 1
$ws1.Cells.Item(19,4).Value = 269.25
$ws1.Cells.Item(19,5).Value = 269.25
$ws1.Cells.Item(19,6).Value = "🟡 Observer"
$ws1.Cells.Item(19,7).Value = "➖ Neutre"

# Row 20: BRVM - INDUSTRIE                (**)
$ws1.Cells.Item(20,1).Value = "BRVM - INDUSTRIE                (**)"
$ws1.Cells.Item(20,2).Value = 0
$ws1.Cells.Item(20,3).Value = 1
$ws1.Cells.Item(20,4).Value = 269.1
$ws1.Cells.Item(20,5).Value = 269.1
$ws1.Cells.Item(20,6).Value = "🟡 Observer"
$ws1.Cells.Item(20,7).Value = "➖ Neutre"

# Row 21: BRVM - INDUSTRIE                  (**)
$ws1.Cells.Item(21,1).Value = "BRVM - INDUSTRIE                  (**)"
$ws1.Cells.Item(21,2).Value = 0
$ws1.Cells.Item(21,3).Value = 1
$ws1.Cells.Item(21,4).Value = 266.59
$ws1.Cells.Item(21,5).Value = 266.59
$ws1.Cells.Item(21,6).Value = "🟡 Observer"
$ws1.Cells.Item(21,7).Value = "➖ Neutre"

# Row 22: BRVM-PRINCIPAL                   (**)
$ws1.Cells.Item(22,1).Value = "BRVM-PRINCIPAL                   (**)"
$ws1.Cells.Item(22,2).Value = 0
$ws1.Cells.Item(22,3).Value = 1
$ws1.Cells.Item(22,4).Value = 222.15
$ws1.Cells.Item(22,5).Value = 222.15
$ws1.Cells.Item(22,6).Value = "🟡 Observer"
$ws1.Cells.Item(22,7).Value = "➖ Neutre"

# Row 23: SETAO CI (STAC)
$ws1.Cells.Item(23,1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(23,2).Value = 2
$ws1.Cells.Item(23,3).Value = 0
$ws1.Cells.Item(23,4).Value = 7.84
$ws1.Cells.Item(23,5).Value = 5.02
$ws1.Cells.Item(23,6).Value = "🟡 Observer"
$ws1.Cells.Item(23,7).Value = "➖ Neutre"

# Row 24: FILTISAC CI (FTSC)
$ws1.Cells.Item(24,1).Value = "FILTISAC CI (FTSC)"
$ws1.Cells.Item(24,2).Value = 1
$ws1.Cells.Item(24,3).Value = 0
$ws1.Cells.Item(24,4).Value = 7.43
$ws1.Cells.Item(24,5).Value = 7.43
$ws1.Cells.Item(24,6).Value = "🟡 Observer"
$ws1.Cells.Item(24,7).Value = "➖ Neutre"

# Row 25: ECOBANK COTE D''IVOIRE (ECOC)
$ws1.Cells.Item(25,1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(25,2).Value = 2
$ws1.Cells.Item(25,3).Value = 1
$ws1.Cells.Item(25,4).Value = 6.23
$ws1.Cells.Item(25,5).Value = 4.99
$ws1.Cells.Item(25,6).Value = "🟡 Observer"
$ws1.Cells.Item(25,7).Value = "👀 À surveiller"

# Row 26: SOCIETE GENERALE COTE D'IVOIRE (SGBC)
$ws1.Cells.Item(26,1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Cells.Item(26,2).Value = 2
$ws1.Cells.Item(26,3).Value = 0
$ws1.Cells.Item(26,4).Value = 5.6
$ws1.Cells.Item(26,5).Value = 1.45
$ws1.Cells.Item(26,6).Value = "🟡 Observer"
$ws1.Cells.Item(26,7).Value = "➖ Neutre"

# Row 27: ORAGROUP TOGO (ORGT)
$ws1.Cells.Item(27,1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(27,2).Value = 1
$ws1.Cells.Item(27,3).Value = 1
$ws1.Cells.Item(27,4).Value = 3.38
$ws1.Cells.Item(27,5).Value = -2.53
$ws1.Cells.Item(27,6).Value = "🟡 Observer"
$ws1.Cells.Item(27,7).Value = "👀 À surveiller"

# Row 28: CFAO MOTORS CI (CFAC)
$ws1.Cells.Item(28,1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(28,2).Value = 1
$ws1.Cells.Item(28,3).Value = 1
$ws1.Cells.Item(28,4).Value = 3.31
$ws1.Cells.Item(28,5).Value = -1.4
$ws1.Cells.Item(28,6).Value = "🟡 Observer"
$ws1.Cells.Item(28,7).Value = "👀 À surveiller"

# Row 29: SMB CI (SMBC)
$ws1.Cells.Item(29,1).Value = "SMB CI (SMBC)"
$ws1.Cells.Item(29,2).Value = 1
$ws1.Cells.Item(29,3).Value = 0
$ws1.Cells.Item(29,4).Value = 3.05
$ws1.Cells.Item(29,5).Value = 3.05
$ws1.Cells.Item(29,6).Value = "🟡 Observer"
$ws1.Cells.Item(29,7).Value = "➖ Neutre"

# Row 30: BANK OF AFRICA BN (BOAB)
$ws1.Cells.Item(30,1).Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Cells.Item(30,2).Value = 1
$ws1.Cells.Item(30,3).Value = 0
$ws1.Cells.Item(30,4).Value = 2.9
$ws1.Cells.Item(30,5).Value = 2.9
$ws1.Cells.Item(30,6).Value = "🟡 Observer"
$ws1.Cells.Item(30,7).Value = "➖ Neutre"

# Row 31: NEI-CEDA CI (NEIC)
$ws1.Cells.Item(31,1).Value = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(31,2).Value = 1
$ws1.Cells.Item(31,3).Value = 1
$ws1.Cells.Item(31,4).Value = 2.81
$ws1.Cells.Item(31,5).Value = 4.97
$ws1.Cells.Item(31,6).Value = "🟡 Observer"
$ws1.Cells.Item(31,7).Value = "👀 À surveiller"

# Row 32: SICABLE CI (CABC)
$ws1.Cells.Item(32,1).Value = "SICABLE CI (CABC)"
$ws1.Cells.Item(32,2).Value = 2
$ws1.Cells.Item(32,3).Value = 1
$ws1.Cells.Item(32,4).Value = 1.2
$ws1.Cells.Item(32,5).Value = 3
$ws1.Cells.Item(32,6).Value = "🟡 Observer"
$ws1.Cells.Item(32,7).Value = "👀 À surveiller"

# Row 33: BANK OF AFRICA SENEGAL (BOAS)
$ws1.Cells.Item(33,1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Cells.Item(33,2).Value = 1
$ws1.Cells.Item(33,3).Value = 0
$ws1.Cells.Item(33,4).Value = 0.82
$ws1.Cells.Item(33,5).Value = 0.82
$ws1.Cells.Item(33,6).Value = "🟡 Observer"
$ws1.Cells.Item(33,7).Value = "➖ Neutre"

# Row 34: TOTALENERGIES MARKETING SN (TTLS)
$ws1.Cells.Item(34,1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Cells.Item(34,2).Value = 1
$ws1.Cells.Item(34,3).Value = 1
$ws1.Cells.Item(34,4).Value = 0.67
$ws1.Cells.Item(34,5).Value = -2.35
$ws1.Cells.Item(34,6).Value = "🟡 Observer"
$ws1.Cells.Item(34,7).Value = "👀 À surveiller"

# Row 35: SERVAIR ABIDJAN CI (ABJC)
$ws1.Cells.Item(35,1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Cells.Item(35,2).Value = 1
$ws1.Cells.Item(35,3).Value = 0
$ws1.Cells.Item(35,4).Value = 0.6
$ws1.Cells.Item(35,5).Value = 0.6
$ws1.Cells.Item(35,6).Value = "🟡 Observer"
$ws1.Cells.Item(35,7).Value = "➖ Neutre"

# Row 36: TOTAL
$ws1.Cells.Item(36,1).Value = "TOTAL"
$ws1.Cells.Item(36,2).Value = 0
$ws1.Cells.Item(36,3).Value = 3
$ws1.Cells.Item(36,4).Value = 0
$ws1.Cells.Item(36,5).Value = 0
$ws1.Cells.Item(36,6).Value = "🟡 Observer"
$ws1.Cells.Item(36,7).Value = "➖ Neutre"

# Row 37: BERNABE CI (BNBC)
$ws1.Cells.Item(37,1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(37,2).Value = 1
$ws1.Cells.Item(37,3).Value = 1
$ws1.Cells.Item(37,4).Value = -0.32
$ws1.Cells.Item(37,5).Value = -2.52
$ws1.Cells.Item(37,6).Value = "🟡 Observer"
$ws1.Cells.Item(37,7).Value = "👀 À surveiller"

# Row 38: ORANGE COTE D'IVOIRE (ORAC)
$ws1.Cells.Item(38,1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Cells.Item(38,2).Value = 1
$ws1.Cells.Item(38,3).Value = 2
$ws1.Cells.Item(38,4).Value = -0.94
$ws1.Cells.Item(38,5).Value = -4.18
$ws1.Cells.Item(38,6).Value = "🟡 Observer"
$ws1.Cells.Item(38,7).Value = "👀 À surveiller"

# Row 39: SUCRIVOIRE (SCRC)
$ws1.Cells.Item(39,1).Value = "SUCRIVOIRE (SCRC)"
$ws1.Cells.Item(39,2).Value = 0
$ws1.Cells.Item(39,3).Value = 1
$ws1.Cells.Item(39,4).Value = -1
$ws1.Cells.Item(39,5).Value = -1
$ws1.Cells.Item(39,6).Value = "🟡 Observer"
$ws1.Cells.Item(39,7).Value = "➖ Neutre"

# Row 40: SOLIBRA CI (SLBC)
$ws1.Cells.Item(40,1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(40,2).Value = 0
$ws1.Cells.Item(40,3).Value = 1
$ws1.Cells.Item(40,4).Value = -1.45
$ws1.Cells.Item(40,5).Value = -1.45
$ws1.Cells.Item(40,6).Value = "🟡 Observer"
$ws1.Cells.Item(40,7).Value = "➖ Neutre"

# Row 41: AIR LIQUIDE CI (SIVC)
$ws1.Cells.Item(41,1).Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Cells.Item(41,2).Value = 1
$ws1.Cells.Item(41,3).Value = 1
$ws1.Cells.Item(41,4).Value = -1.73
$ws1.Cells.Item(41,5).Value = 5.26
$ws1.Cells.Item(41,6).Value = "🟡 Observer"
$ws1.Cells.Item(41,7).Value = "👀 À surveiller"

# Row 42: SAPH CI (SPHC)
$ws1.Cells.Item(42,1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(42,2).Value = 0
$ws1.Cells.Item(42,3).Value = 1
$ws1.Cells.Item(42,4).Value = -1.84
$ws1.Cells.Item(42,5).Value = -1.84
$ws1.Cells.Item(42,6).Value = "🟡 Observer"
$ws1.Cells.Item(42,7).Value = "➖ Neutre"

# Row 43: UNIWAX CI (UNXC)
$ws1.Cells.Item(43,1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(43,2).Value = 0
$ws1.Cells.Item(43,3).Value = 1
$ws1.Cells.Item(43,4).Value = -2.68
$ws1.Cells.Item(43,5).Value = -2.68
$ws1.Cells.Item(43,6).Value = "🟡 Observer"
$ws1.Cells.Item(43,7).Value = "➖ Neutre"

# Row 44: VIVO ENERGY CI (SHEC)
$ws1.Cells.Item(44,1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(44,2).Value = 0
$ws1.Cells.Item(44,3).Value = 2
$ws1.Cells.Item(44,4).Value = -3.56
$ws1.Cells.Item(44,5).Value = -1.94
$ws1.Cells.Item(44,6).Value = "🟡 Observer"
$ws1.Cells.Item(44,7).Value = "➖ Neutre"

# Row 45: ECOBANK TRANS. INCORP. TG (ETIT)
$ws1.Cells.Item(45,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(45,2).Value = 0
$ws1.Cells.Item(45,3).Value = 1
$ws1.Cells.Item(45,4).Value = -4.35
$ws1.Cells.Item(45,5).Value = -4.35
$ws1.Cells.Item(45,6).Value = "🟡 Observer"
$ws1.Cells.Item(45,7).Value = "➖ Neutre"

# Row 46: BICI CI (BICC)
$ws1.Cells.Item(46,1).Value = "BICI CI (BICC)"
$ws1.Cells.Item(46,2).Value = 0
$ws1.Cells.Item(46,3).Value = 2
$ws1.Cells.Item(46,4).Value = -5.44
$ws1.Cells.Item(46,5).Value = -2.21
$ws1.Cells.Item(46,6).Value = "🟡 Observer"
$ws1.Cells.Item(46,7).Value = "➖ Neutre"

# Row 47: SICOR CI (SICC)
$ws1.Cells.Item(47,1).Value = "SICOR CI (SICC)"
$ws1.Cells.Item(47,2).Value = 0
$ws1.Cells.Item(47,3).Value = 1
$ws1.Cells.Item(47,4).Value = -6.87
$ws1.Cells.Item(47,5).Value = -6.87
$ws1.Cells.Item(47,6).Value = "🟡 Observer"
$ws1.Cells.Item(47,7).Value = "➖ Neutre"

# --- Sheet "Top_YTD": refresh YTD progression ranking (rows 2-11) ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 2: BRVM - SERVICES PUBLICS
$ws2.Cells.Item(2,1).Value = "BRVM - SERVICES PUBLICS"
$ws2.Cells.Item(2,2).Value = 10271314.08

# Row 3: NEI-CEDA CI
$ws2.Cells.Item(3,1).Value = "NEI-CEDA CI"
$ws2.Cells.Item(3,2).Value = 1081631.75

# Row 4: AIR LIQUIDE CI
$ws2.Cells.Item(4,1).Value = "AIR LIQUIDE CI"
$ws2.Cells.Item(4,2).Value = 383807.6

# Row 5: BRVM - AUTRES SECTEURS
$ws2.Cells.Item(5,1).Value = "BRVM - AUTRES SECTEURS"
$ws2.Cells.Item(5,2).Value = 227718.18

# Row 6: BRVM - DISTRIBUTION
$ws2.Cells.Item(6,1).Value = "BRVM - DISTRIBUTION"
$ws2.Cells.Item(6,2).Value = 167175.51

# Row 7: BRVM - TRANSPORT
$ws2.Cells.Item(7,1).Value = "BRVM - TRANSPORT"
$ws2.Cells.Item(7,2).Value = 43304.75

# Row 8: BRVM - AGRICULTURE
$ws2.Cells.Item(8,1).Value = "BRVM - AGRICULTURE"
$ws2.Cells.Item(8,2).Value = 38358.53

# Row 9: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws2.Cells.Item(9,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Cells.Item(9,2).Value = 6694.97

# Row 10: BRVM - FINANCES
$ws2.Cells.Item(10,1).Value = "BRVM - FINANCES"
$ws2.Cells.Item(10,2).Value = 3631.85

# Row 11: BRVM - SERVICES FINANCIERS
$ws2.Cells.Item(11,1).Value = "BRVM - SERVICES FINANCIERS"
$ws2.Cells.Item(11,2).Value = 3481.09
